$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1329
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 1391.9
$ws.Range("K46").Value = 2100
$ws.Range("L46").Value = 4175.700000000001
$ws.Range("M46").Value = -1981
$ws.Range("N46").Value = -4413.700000000001
$ws.Range("H51").Value = 8011624
$ws.Range("I51").Value = 33367984
$ws.Range("J51").Value = 4352.6313
$ws.Range("K51").Value = 33367984
$ws.Range("L51").Value = 4352.6313
$ws.Range("M51").Value = -33367500
$ws.Range("N51").Value = -5320.6313
$ws.Range("H60").Value = 1329
$ws.Range("I60").Value = 700
$ws.Range("J60").Value = 1391.9
$ws.Range("K60").Value = 2100
$ws.Range("L60").Value = 4175.700000000001
$ws.Range("M60").Value = -1616
$ws.Range("N60").Value = -5143.700000000001
$ws.Range("H99").Value = 265.0476
$ws.Range("I99").Value = 271.8421
$ws.Range("J99").Value = 200.5
$ws.Range("K99").Value = 815.5263
$ws.Range("L99").Value = 601.5
$ws.Range("M99").Value = 682.4737
$ws.Range("N99").Value = -3597.5
$ws.Range("H118").Value = 631.6667
$ws.Range("I118").Value = 447.5
$ws.Range("K118").Value = 1342.5
$ws.Range("M118").Value = 314.5
$ws.Range("H129").Value = 76924620
$ws.Range("I129").Value = 142858420
$ws.Range("J129").Value = 1863.3334
$ws.Range("K129").Value = 428575260
$ws.Range("L129").Value = 5590.0002
$ws.Range("M129").Value = -428570260
$ws.Range("N129").Value = -15590.0002
$ws.Range("H138").Value = 2950.798
$ws.Range("I138").Value = 1164.5151
$ws.Range("J138").Value = 3843.9395
$ws.Range("K138").Value = 3493.5453
$ws.Range("L138").Value = 11531.8185
$ws.Range("M138").Value = 1646.4547
$ws.Range("N138").Value = -21811.8185
$ws.Range("H141").Value = 7400.0835
$ws.Range("I141").Value = 3176.7778
$ws.Range("J141").Value = 20070
$ws.Range("K141").Value = 9530.3334
$ws.Range("L141").Value = 60210
$ws.Range("M141").Value = -4350.3334
$ws.Range("N141").Value = -70570

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26338688
$ws.Range("I32").Value = 45469824
$ws.Range("J32").Value = 33374.938
$ws.Range("K32").Value = 45469824
$ws.Range("L32").Value = 33374.938
$ws.Range("M32").Value = -45469537
$ws.Range("N32").Value = -33948.938
$ws.Range("H74").Value = 2182.805
$ws.Range("I74").Value = 2344.8333
$ws.Range("J74").Value = 1740.909
$ws.Range("K74").Value = 2344.8333
$ws.Range("L74").Value = 1740.909
$ws.Range("M74").Value = -1470.8333
$ws.Range("N74").Value = -3488.909
$ws.Range("H77").Value = 2182.805
$ws.Range("I77").Value = 2344.8333
$ws.Range("J77").Value = 1740.909
$ws.Range("K77").Value = 11724.1665
$ws.Range("L77").Value = 8704.545
$ws.Range("M77").Value = -7356.166499999999
$ws.Range("N77").Value = -17440.545
$ws.Range("H110").Value = 1276.5
$ws.Range("I110").Value = 702
$ws.Range("J110").Value = 2343.4285
$ws.Range("K110").Value = 702
$ws.Range("L110").Value = 2343.4285
$ws.Range("M110").Value = 1343
$ws.Range("N110").Value = -6433.4285
$ws.Range("H132").Value = 30307470
$ws.Range("I132").Value = 50001160
$ws.Range("K132").Value = 150003480
$ws.Range("M132").Value = -150000950

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 13933.333
$ws.Range("J44").Value = 13933.333
$ws.Range("L44").Value = 13933.333
$ws.Range("N44").Value = -14927.333
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3667.5745
$ws.Range("I31").Value = 2101.5
$ws.Range("J31").Value = 5301.7393
$ws.Range("K31").Value = 2101.5
$ws.Range("L31").Value = 5301.7393
$ws.Range("M31").Value = -1806.5
$ws.Range("N31").Value = -5891.7393
$ws.Range("H34").Value = 3667.5745
$ws.Range("I34").Value = 2101.5
$ws.Range("J34").Value = 5301.7393
$ws.Range("K34").Value = 2101.5
$ws.Range("L34").Value = 5301.7393
$ws.Range("M34").Value = -1899.5
$ws.Range("N34").Value = -5705.7393
$ws.Range("H132").Value = 9808255
$ws.Range("I132").Value = 977.2917
$ws.Range("K132").Value = 2931.8751
$ws.Range("M132").Value = -401.8751000000002

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 16669324
$ws.Range("J22").Value = 3188.4
$ws.Range("L22").Value = 9565.200000000001
$ws.Range("N22").Value = -9903.200000000001
$ws.Range("H27").Value = 16669324
$ws.Range("J27").Value = 3188.4
$ws.Range("L27").Value = 9565.200000000001
$ws.Range("N27").Value = -9769.200000000001
$ws.Range("H113").Value = 45833784
$ws.Range("I113").Value = 41667068
$ws.Range("J113").Value = 50000504
$ws.Range("K113").Value = 125001204
$ws.Range("L113").Value = 150001512
$ws.Range("M113").Value = -124999034
$ws.Range("N113").Value = -150005852
$ws.Range("H131").Value = 777.45
$ws.Range("I131").Value = 491.66666
$ws.Range("J131").Value = 795.69147
$ws.Range("K131").Value = 1474.99998
$ws.Range("L131").Value = 2387.07441
$ws.Range("M131").Value = 3565.00002
$ws.Range("N131").Value = -12467.07441
$ws.Range("H132").Value = 22733364
$ws.Range("I132").Value = 1057
$ws.Range("J132").Value = 33341774
$ws.Range("K132").Value = 9513
$ws.Range("L132").Value = 300075966
$ws.Range("M132").Value = -6983
$ws.Range("N132").Value = -300081026

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2082.2
$ws.Range("I113").Value = 2803.6667
$ws.Range("K113").Value = 2803.6667
$ws.Range("M113").Value = -633.6667000000002
$ws.Range("H116").Value = 20586.666
$ws.Range("J116").Value = 20586.666
$ws.Range("L116").Value = 20586.666
$ws.Range("N116").Value = -29764.666
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 34950
$ws.Range("J120").Value = 34950
$ws.Range("L120").Value = 34950
$ws.Range("N120").Value = -44626

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2534165.5
$ws.Range("I22").Value = 12658227
$ws.Range("J22").Value = 3150
$ws.Range("K22").Value = 12658227
$ws.Range("L22").Value = 3150
$ws.Range("M22").Value = -12657932
$ws.Range("N22").Value = -3740
$ws.Range("H27").Value = 2534165.5
$ws.Range("I27").Value = 12658227
$ws.Range("J27").Value = 3150
$ws.Range("K27").Value = 12658227
$ws.Range("L27").Value = 3150
$ws.Range("M27").Value = -12658120
$ws.Range("N27").Value = -3364
$ws.Range("H46").Value = 8333992.5
$ws.Range("I46").Value = 20833732
$ws.Range("J46").Value = 833.3333
$ws.Range("K46").Value = 20833732
$ws.Range("L46").Value = 833.3333
$ws.Range("M46").Value = -20833544
$ws.Range("N46").Value = -1209.3333
$ws.Range("H55").Value = 111122350
$ws.Range("I55").Value = 33434
$ws.Range("K55").Value = 33434
$ws.Range("M55").Value = -33261
$ws.Range("H61").Value = 1641.5385
$ws.Range("I61").Value = 1294.5
$ws.Range("J61").Value = 2798.3333
$ws.Range("K61").Value = 1294.5
$ws.Range("L61").Value = 2798.3333
$ws.Range("M61").Value = -1092.5
$ws.Range("N61").Value = -3202.3333
$ws.Range("H113").Value = 1641.5385
$ws.Range("I113").Value = 1294.5
$ws.Range("J113").Value = 2798.3333
$ws.Range("K113").Value = 1294.5
$ws.Range("L113").Value = 2798.3333
$ws.Range("M113").Value = 875.5
$ws.Range("N113").Value = -7138.3333
$ws.Range("H136").Value = 3693.361
$ws.Range("I136").Value = 3132.074
$ws.Range("J136").Value = 5377.222
$ws.Range("K136").Value = 9396.222
$ws.Range("L136").Value = 16131.666
$ws.Range("M136").Value = -6846.222
$ws.Range("N136").Value = -21231.666

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").ClearContents()
